$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Factor buckets moved from quintiles (0-4, 5 columns per market) to
# quartiles (0-3, 4 columns per market) -> columns J and K are no longer
# needed at all, so remove them outright (shifts dimension/spans too).
$ws.Columns("J:K").Delete()

# --- Re-merge the header band: was B1:F1 / G1:K1 (5 cols each),
#     now B1:E1 / F1:I1 (4 cols each) ---
$ws.Range("B1:I1").UnMerge()
$ws.Range("B1:E1").Merge()
$ws.Range("F1:I1").Merge()

# Re-apply the original header formatting (bold + border + centered) to
# every cell in the header row, since re-merging recomputes per-edge
# border styles on the affected cells.
$ws.Range("A1").Copy()
$ws.Range("A1:I1").PasteSpecial(-4122)

# "Emerging" band now starts at F1 (used to be G1)
$ws.Range("F1").Value = "Emerging"

# --- Row 2 (quartile index row): 0,1,2,3 under each of the two merged
#     bands instead of 0,1,2,3,4 ---
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 3

# --- Row 4 (forward-return data): new quartile-based values ---
$ws.Range("B4").Value = 0.008237511945334871
$ws.Range("C4").Value = 0.005741652032354031
$ws.Range("D4").Value = 0.007357543898553393
$ws.Range("E4").Value = 0.009645294618850515
$ws.Range("F4").Value = 0.01668172215606266
$ws.Range("G4").Value = 0.006111143533269399
$ws.Range("H4").Value = 0.0105130582725587
$ws.Range("I4").Value = 0.01271031198261327
